$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# The sheet is a flat "Conta/Nome/Saldo" export. A new record
# (005232019 / PEDRO / 3000) needs to land right before the existing
# "004392159 / RODRIGO / 900.21" row, which is spreadsheet row 15
# (row 1 is the header). Insert a blank row there first so every
# following row shifts down by one, then fill it in.
$ws.Rows(15).Insert()

# "Conta" values are account numbers with leading zeros, stored as text
# in the source file (inlineStr), so force text formatting before
# assigning the numeric-looking string - otherwise Excel would coerce
# "005232019" into the number 5232019 and drop the leading zeros.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "005232019"
$ws.Range("B15").Value = "PEDRO"
$ws.Range("C15").Value = 3000
